# Auto-generated script to update cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.143.13"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "'  -0.40%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'1.896.41"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = "'  -0.57%  "
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'  +0.29%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'307.48"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = "'  +0.02%  "
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'  +0.21%  "
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'  -0.89%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.3768"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'  -0.36%  "
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "'  +0.35%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'21.19"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'  -0.34%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.9008"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'  +0.05%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.08212"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'  -0.86%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'1.943.06"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'96.15"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'  +0.77%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'5.331"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'  +0.96%  "
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'  +0.28%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'0.000008623"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'  +0.22%  "
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'  +0.89%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'1.003"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "'  +0.38%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'27.172.12"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'5.097"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'  +0.64%  "
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'  +0.61%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'6.428"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'  -0.47%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'149.08"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'  +2.14%  "
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'  +0.45%  "
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'  +0.18%  "
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'  +0.02%  "
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "'  +0.51%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'4.807"
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "'  -0.16%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'4.866"
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.Value = "'  -2.17%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'0.09214"
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = "'  +0.16%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'0.7968"
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "'  -0.87%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'0.05035"
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "'  -0.99%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'1.218"
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = "'  -1.65%  "
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = "'  +2.24%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'2.956"
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "'  +0.69%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'2.613"
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.Value = "'  +1.60%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'0.5731"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = "'  -0.19%  "
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'  +0.33%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'9.008"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'  -0.68%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'6.561"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "'  -1.17%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'116.27"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "'  -1.84%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.1514"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'  -0.23%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.4873"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'  +0.62%  "
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'  +0.26%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'10.11"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'  -0.38%  "
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'  +0.33%  "
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'  +2.16%  "
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'63.74"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'  +0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.05931"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'  -0.33%  "
$c.Style = "Normal"
